# Update the "想去人数" (want-to-go count, column F) figures that changed
# between the two generated-data snapshots, on both the "展览" (sheet 1)
# and "全部类型" (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 100
$wsExhibit.Range("F4").Value  = 7317
$wsExhibit.Range("F5").Value  = 272
$wsExhibit.Range("F6").Value  = 428
$wsExhibit.Range("F7").Value  = 3831
$wsExhibit.Range("F8").Value  = 309
$wsExhibit.Range("F9").Value  = 542
$wsExhibit.Range("F11").Value = 616
$wsExhibit.Range("F12").Value = 106

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value  = 100
$wsAll.Range("F5").Value  = 7317
$wsAll.Range("F7").Value  = 272
$wsAll.Range("F8").Value  = 428
$wsAll.Range("F9").Value  = 3831
$wsAll.Range("F10").Value = 309
$wsAll.Range("F11").Value = 542
$wsAll.Range("F13").Value = 616
$wsAll.Range("F14").Value = 106
